$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format for numeric-looking values so Excel does not
# reinterpret them as numbers/dates and strip formatting (trailing zeros etc).
$textCells = @("D5", "D6", "D7", "D8", "D11", "D15", "D16", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.036.85"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.503.18"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "534.42"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "136.36"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D11").Value = "5.39"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "2.945.73"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "58.935.67"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "22.73"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "0.0000139"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "2.507.26"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "4.26"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "323.60"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "6.08"
$ws.Range("E21").Value = "  +4.39%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "65.22"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").Value = "0.422"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "0.0₃0764"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "6.50"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").Value = "170.20"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "1.75"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "18.37"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "4.06"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "0.802"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "282.81"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").Value = "130.51"
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "10.92"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "0.0925"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "0.0501"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "17.33"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "1.757.19"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  +0.20%  "
